# Add 2022-Q4 data
#
# Before:  总计 | 2022-Q3 | 2022-Q2
# After:   总计 | 2022-Q4 | 2022-Q3 | 2022-Q2
#
# The new "2022-Q4" sheet is a sibling of the existing quarter sheets
# (same columns/layout), so we clone the "2022-Q3" sheet (to inherit its
# formatting/styles exactly) and then overwrite its data with the Q4
# figures. The "总计" (totals) sheet gets a new top data row for Q4 and
# its existing rows shift down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" summary sheet: insert the 2022-Q4 totals row, push the rest
#    down by one (values only - same row indices the diff touches).
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Old row 3 (2022-Q2) slides down to row 4, keeping its values.
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q2"
$totals.Range("C4").Value = 9
$totals.Range("D4").Value = 1.46
$totals.Range("A3").Copy()
$totals.Range("A4").PasteSpecial(-4122)

# Old row 2 (2022-Q3) slides down to row 3, keeping its values.
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("D3").Value = 1.2

# Row 2 becomes the brand-new 2022-Q4 totals.
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("D2").Value = 1.24

# ---------------------------------------------------------------------
# 2. Clone the "2022-Q3" sheet (current position 2) into a new sheet
#    placed right after "总计", then rename it "2022-Q4".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $totals)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 3. Overwrite the cloned sheet's figures with the real 2022-Q4 numbers.
#    Columns D/E/F/G are text cells (quote-prefixed so they stay text,
#    matching the source data which stores them as strings); column H
#    is numeric.
# ---------------------------------------------------------------------
$q4rows = @(
    @(2,  "3.23", "83.79", "7.87", "0.2542", 1),
    @(3,  "3.25", "94.43", "7.78", "0.2528", 2),
    @(4,  "3.24", "68.93", "7.33", "0.2375", 5),
    @(5,  "3.06", "94.24", "5.21", "0.1594", 7),
    @(6,  "4.05", "90.08", "3.92", "0.1588", 7),
    @(7,  "0.81", "83.79", "7.87", "0.0637", 1),
    @(8,  "0.77", "68.93", "7.33", "0.0564", 5),
    @(9,  "0.45", "94.43", "7.78", "0.0350", 2),
    @(10, "0.48", "94.24", "5.21", "0.0250", 7)
)

foreach ($row in $q4rows) {
    $r = $row[0]
    $q4.Cells.Item($r, 4).Value = "'" + $row[1]
    $q4.Cells.Item($r, 5).Value = "'" + $row[2]
    $q4.Cells.Item($r, 6).Value = "'" + $row[3]
    $q4.Cells.Item($r, 7).Value = "'" + $row[4]
    $q4.Cells.Item($r, 8).Value = $row[5]
}

# Restore the originally-active tab (the copy operation above would
# otherwise leave the new 2022-Q4 sheet selected).
$wb.Worksheets.Item("2022-Q2").Activate()
